# AWS bash cmd line
# Refresh the "last status check" timestamp and replace the Shell
# Olomoucka row's delta/timestamp columns (previously stored as literal
# text) with real numeric values, matching the other rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: bump the status-check time from 09:15 to 09:30.
$ws.Range("F1").Value = "Last status check on: 26.01.2022 09:30"

# Row 6 (Shell Olomoucka): D6 was the text "+0.2" -> numeric 0.2.
$ws.Range("D6").Value = 0.2

# Row 6: E6 was the text "2022-01-26 09:15:12" -> numeric datetime
# serial (matches the date/time format already used by sibling rows).
$ws.Range("E6").Value = 44587.38555555556
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
